# "Corr/total marks" — update the marksheet's correct/total mark figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer marking value
$ws.Range("B11").Value = 5

# Total row: total marks obtained, and the "correct/total" summary label
$ws.Range("B12").Value = 70
$ws.Range("E12").Value = "70/140"
